$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F is "dSF" (delta S final); this edit repulls/recomputes it for
# several rows (and a couple of related columns on row 8) per the
# "repull data, push all data, mean calculation" commit.

$ws.Range("F5").Value = -9
$ws.Range("F7").Value = -8
$ws.Range("E8").Value = 4
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 7
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = -2
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("F19").Value = -3
$ws.Range("F22").Value = 0
$ws.Range("F25").Value = -3
$ws.Range("F26").Value = -4
$ws.Range("F29").Value = -9
$ws.Range("F32").Value = -2
$ws.Range("F33").Value = -4
$ws.Range("F34").Value = -4
$ws.Range("F35").Value = 8
$ws.Range("F37").Value = -3
$ws.Range("F39").Value = 0
$ws.Range("F41").Value = -5
$ws.Range("F44").Value = 6
$ws.Range("F47").Value = -3
$ws.Range("F49").Value = -2
$ws.Range("F51").Value = -10
$ws.Range("F60").Value = -3
$ws.Range("F61").Value = 6
$ws.Range("F63").Value = -2
$ws.Range("F71").Value = 2
$ws.Range("F72").Value = 6
